$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 38575
$ws.Range("D2").Value = 55775598
$ws.Range("C3").Value = 92571
$ws.Range("D3").Value = 135689428
$ws.Range("C4").Value = 31651
$ws.Range("D4").Value = 46871320
$ws.Range("C5").Value = 8876
$ws.Range("D5").Value = 13191397
$ws.Range("C6").Value = 2056
$ws.Range("D6").Value = 3055971
$ws.Range("C12").Value = 42032
$ws.Range("D12").Value = 57003801
$ws.Range("C13").Value = 9847
$ws.Range("D13").Value = 14242628
$ws.Range("C14").Value = 26354
$ws.Range("D14").Value = 38639335
$ws.Range("C15").Value = 8425
$ws.Range("D15").Value = 12502978
$ws.Range("C16").Value = 2192
$ws.Range("D16").Value = 3257930
$ws.Range("C17").Value = 428
$ws.Range("D17").Value = 631123
$ws.Range("C20").Value = 10374
$ws.Range("D20").Value = 13726087
$ws.Range("C21").Value = 13626
$ws.Range("D21").Value = 19665766
$ws.Range("C22").Value = 32112
$ws.Range("D22").Value = 47119323
$ws.Range("C23").Value = 10365
$ws.Range("D23").Value = 15406210
$ws.Range("C24").Value = 2678
$ws.Range("D24").Value = 3981771
$ws.Range("C27").Value = 11879
$ws.Range("D27").Value = 15856612
$ws.Range("C28").Value = 7834
$ws.Range("D28").Value = 11336407
$ws.Range("C29").Value = 22900
$ws.Range("D29").Value = 33613315
$ws.Range("C30").Value = 7911
$ws.Range("D30").Value = 11766591
$ws.Range("C31").Value = 2001
$ws.Range("D31").Value = 2985919
$ws.Range("C32").Value = 376
$ws.Range("D32").Value = 561415
$ws.Range("C34").Value = 8444
$ws.Range("D34").Value = 11154320
$ws.Range("C35").Value = 3335
$ws.Range("D35").Value = 4817652
$ws.Range("C36").Value = 7988
$ws.Range("D36").Value = 11665821
$ws.Range("C37").Value = 3223
$ws.Range("D37").Value = 4777461
$ws.Range("C41").Value = 2521
$ws.Range("D41").Value = 3406222
$ws.Range("C42").Value = 17618
$ws.Range("D42").Value = 25476710
$ws.Range("C43").Value = 51994
$ws.Range("D43").Value = 76209745
$ws.Range("C44").Value = 19254
$ws.Range("D44").Value = 28595379
$ws.Range("C45").Value = 5706
$ws.Range("D45").Value = 8494260
$ws.Range("C46").Value = 1241
$ws.Range("D46").Value = 1852045
$ws.Range("C50").Value = 17047
$ws.Range("D50").Value = 22656252
$ws.Range("C51").Value = 2114
$ws.Range("D51").Value = 3066878
$ws.Range("C52").Value = 7156
$ws.Range("D52").Value = 10515701
$ws.Range("C53").Value = 2410
$ws.Range("D53").Value = 3599464
$ws.Range("C55").Value = 196
$ws.Range("D55").Value = 290226
$ws.Range("C57").Value = 7270
$ws.Range("D57").Value = 9995619
$ws.Range("C58").Value = 1130
$ws.Range("D58").Value = 1881544
$ws.Range("C59").Value = 2768
$ws.Range("D59").Value = 4587488
$ws.Range("C60").Value = 1089
$ws.Range("D60").Value = 1805338
$ws.Range("C61").Value = 374
$ws.Range("D61").Value = 623383
$ws.Range("C64").Value = 1650
$ws.Range("D64").Value = 2544202
$ws.Range("C65").Value = 15696
$ws.Range("D65").Value = 22669283
$ws.Range("C66").Value = 45447
$ws.Range("D66").Value = 66492807
$ws.Range("C67").Value = 15918
$ws.Range("D67").Value = 23650496
$ws.Range("C68").Value = 4631
$ws.Range("D68").Value = 6897551
$ws.Range("C73").Value = 15332
$ws.Range("D73").Value = 20198546
$ws.Range("C74").Value = 53748
$ws.Range("D74").Value = 78221205
$ws.Range("C75").Value = 150997
$ws.Range("D75").Value = 222445564
$ws.Range("C76").Value = 65274
$ws.Range("D76").Value = 97264950
$ws.Range("C77").Value = 20901
$ws.Range("D77").Value = 31231822
$ws.Range("C78").Value = 4990
$ws.Range("D78").Value = 7453403
$ws.Range("C85").Value = 52997
$ws.Range("D85").Value = 72023669
$ws.Range("C86").Value = 4739
$ws.Range("D86").Value = 6868611
$ws.Range("C87").Value = 11822
$ws.Range("D87").Value = 17364041
$ws.Range("C88").Value = 3944
$ws.Range("D88").Value = 5878583
$ws.Range("C93").Value = 5539
$ws.Range("D93").Value = 7443611
$ws.Range("C94").Value = 1644
$ws.Range("D94").Value = 2368302
$ws.Range("C95").Value = 5303
$ws.Range("D95").Value = 7811878
$ws.Range("C96").Value = 1973
$ws.Range("D96").Value = 2937826
$ws.Range("C101").Value = 3658
$ws.Range("D101").Value = 4840003
$ws.Range("C102").Value = 698
$ws.Range("D102").Value = 1143375
$ws.Range("C107").Value = 11013
$ws.Range("D107").Value = 15975708
$ws.Range("C108").Value = 29634
$ws.Range("D108").Value = 43525745
$ws.Range("C109").Value = 9923
$ws.Range("D109").Value = 14753705
$ws.Range("C110").Value = 2732
$ws.Range("D110").Value = 4073080
$ws.Range("C114").Value = 9958
$ws.Range("D114").Value = 13147631
$ws.Range("C115").Value = 31066
$ws.Range("D115").Value = 44793377
$ws.Range("C116").Value = 67167
$ws.Range("D116").Value = 98283335
$ws.Range("C117").Value = 21657
$ws.Range("D117").Value = 32186038
$ws.Range("C118").Value = 6144
$ws.Range("D118").Value = 9154021
$ws.Range("C119").Value = 1149
$ws.Range("D119").Value = 1717100
$ws.Range("C120").Value = 85
$ws.Range("D120").Value = 123895
$ws.Range("C124").Value = 26240
$ws.Range("D124").Value = 35024911
$ws.Range("C125").Value = 36797
$ws.Range("D125").Value = 53096475
$ws.Range("C126").Value = 78161
$ws.Range("D126").Value = 114282226
$ws.Range("C127").Value = 24194
$ws.Range("D127").Value = 35909284
$ws.Range("C128").Value = 6501
$ws.Range("D128").Value = 9661358
$ws.Range("C129").Value = 1274
$ws.Range("D129").Value = 1894311
$ws.Range("C130").Value = 60
$ws.Range("D130").Value = 88228
$ws.Range("C133").Value = 32346
$ws.Range("D133").Value = 42928398
$ws.Range("C134").Value = 13543
$ws.Range("D134").Value = 19604952
$ws.Range("C135").Value = 32832
$ws.Range("D135").Value = 48215504
$ws.Range("C136").Value = 11633
$ws.Range("D136").Value = 17283587
$ws.Range("C137").Value = 3008
$ws.Range("D137").Value = 4483241
$ws.Range("C141").Value = 10976
$ws.Range("D141").Value = 14629082
$ws.Range("C142").Value = 35920
$ws.Range("D142").Value = 51877876
$ws.Range("C143").Value = 82881
$ws.Range("D143").Value = 121423316
$ws.Range("C144").Value = 24752
$ws.Range("D144").Value = 36772805
$ws.Range("C145").Value = 6501
$ws.Range("D145").Value = 9701067
$ws.Range("C146").Value = 1472
$ws.Range("D146").Value = 2190230
$ws.Range("C149").Value = 29723
$ws.Range("D149").Value = 40066181
